$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Admin delete/edit user functionality implemented -> mark the corresponding
# "Admin Home Screen" / "Admin Approve Ad" / "Admin Reject Ad" / "Admin Edit Ad" /
# "Admin Delete Ad" / "Admin List Users" rows as done ("Yes") in column C,
# keeping the cells' existing styling.
$ws.Range("C34:C39").Value = "Yes"

# Reflect the updated scroll position / active selection from the edit session.
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H4").Select()
